$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44386
$ws.Cells.Item(2, 12).Value = 'Especial'
$ws.Cells.Item(2, 13).Value = 45
$ws.Cells.Item(2, 14).Value = 14000
$ws.Cells.Item(2, 15).Value = 14000
$ws.Cells.Item(2, 16).Value = 14000
$ws.Cells.Item(2, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(2, 19).Value = 933
$ws.Cells.Item(2, 20).Value = 15

# Row 3
$ws.Cells.Item(3, 4).Value = 44386
$ws.Cells.Item(3, 13).Value = 60
$ws.Cells.Item(3, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(3, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(3, 19).Value = 800
$ws.Cells.Item(3, 20).Value = 15

# Row 4
$ws.Cells.Item(4, 4).Value = 44355
$ws.Cells.Item(4, 12).Value = 'Especial'
$ws.Cells.Item(4, 13).Value = 120
$ws.Cells.Item(4, 14).Value = 12000
$ws.Cells.Item(4, 15).Value = 12000
$ws.Cells.Item(4, 16).Value = 12000
$ws.Cells.Item(4, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(4, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 19).Value = 800
$ws.Cells.Item(4, 20).Value = 15

# Row 5
$ws.Cells.Item(5, 4).Value = 44355
$ws.Cells.Item(5, 13).Value = 150
$ws.Cells.Item(5, 14).Value = 10000
$ws.Cells.Item(5, 15).Value = 10000
$ws.Cells.Item(5, 16).Value = 10000
$ws.Cells.Item(5, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(5, 19).Value = 667
$ws.Cells.Item(5, 20).Value = 15

# Row 6
$ws.Cells.Item(6, 4).Value = 44314
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 7
$ws.Cells.Item(6, 14).Value = 230000
$ws.Cells.Item(6, 15).Value = 230000
$ws.Cells.Item(6, 16).Value = 230000
$ws.Cells.Item(6, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(6, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(6, 19).Value = 511
$ws.Cells.Item(6, 20).Value = 450

# Row 7
$ws.Cells.Item(7, 4).Value = 44385
$ws.Cells.Item(7, 12).Value = 'Especial'
$ws.Cells.Item(7, 13).Value = 60
$ws.Cells.Item(7, 14).Value = 14000
$ws.Cells.Item(7, 15).Value = 14000
$ws.Cells.Item(7, 16).Value = 14000
$ws.Cells.Item(7, 19).Value = 933

# Row 8
$ws.Cells.Item(8, 4).Value = 44385
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'

# Row 9
$ws.Cells.Item(9, 4).Value = 44385
$ws.Cells.Item(9, 12).Value = 'Segunda'
$ws.Cells.Item(9, 13).Value = 75
$ws.Cells.Item(9, 18).Value = 'Región de O''Higgins'

# Row 13
$ws.Cells.Item(13, 4).Value = 44376
$ws.Cells.Item(13, 13).Value = 60
$ws.Cells.Item(13, 14).Value = 15000
$ws.Cells.Item(13, 15).Value = 15000
$ws.Cells.Item(13, 16).Value = 15000
$ws.Cells.Item(13, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(13, 19).Value = 1000

# Row 14
$ws.Cells.Item(14, 4).Value = 44376
$ws.Cells.Item(14, 13).Value = 85
$ws.Cells.Item(14, 18).Value = 'Región Metropolitana'

# Row 15
$ws.Cells.Item(15, 4).Value = 44383
$ws.Cells.Item(15, 12).Value = 'Especial'
$ws.Cells.Item(15, 13).Value = 70
$ws.Cells.Item(15, 14).Value = 14000
$ws.Cells.Item(15, 15).Value = 14000
$ws.Cells.Item(15, 16).Value = 14000
$ws.Cells.Item(15, 19).Value = 933

# Row 16
$ws.Cells.Item(16, 4).Value = 44383
$ws.Cells.Item(16, 13).Value = 80
$ws.Cells.Item(16, 14).Value = 12000
$ws.Cells.Item(16, 15).Value = 12000
$ws.Cells.Item(16, 16).Value = 12000
$ws.Cells.Item(16, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(16, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(16, 19).Value = 800
$ws.Cells.Item(16, 20).Value = 15

# Row 17
$ws.Cells.Item(17, 4).Value = 44315
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 5
$ws.Cells.Item(17, 14).Value = 230000
$ws.Cells.Item(17, 15).Value = 230000
$ws.Cells.Item(17, 16).Value = 230000
$ws.Cells.Item(17, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(17, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(17, 19).Value = 511
$ws.Cells.Item(17, 20).Value = 450

# Row 18
$ws.Cells.Item(18, 4).Value = 44307
$ws.Cells.Item(18, 13).Value = 6
$ws.Cells.Item(18, 14).Value = 250000
$ws.Cells.Item(18, 15).Value = 250000
$ws.Cells.Item(18, 16).Value = 250000
$ws.Cells.Item(18, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(18, 19).Value = 556
$ws.Cells.Item(18, 20).Value = 450

# Row 19
$ws.Cells.Item(19, 4).Value = 44364
$ws.Cells.Item(19, 14).Value = 12000
$ws.Cells.Item(19, 15).Value = 12000
$ws.Cells.Item(19, 16).Value = 12000
$ws.Cells.Item(19, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(19, 19).Value = 800

# Row 20
$ws.Cells.Item(20, 4).Value = 44364
$ws.Cells.Item(20, 13).Value = 120
$ws.Cells.Item(20, 14).Value = 10000
$ws.Cells.Item(20, 15).Value = 10000
$ws.Cells.Item(20, 16).Value = 10000
$ws.Cells.Item(20, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(20, 19).Value = 667

# Row 21
$ws.Cells.Item(21, 4).Value = 44364
$ws.Cells.Item(21, 12).Value = 'Segunda'
$ws.Cells.Item(21, 13).Value = 80
$ws.Cells.Item(21, 14).Value = 8000
$ws.Cells.Item(21, 15).Value = 8000
$ws.Cells.Item(21, 16).Value = 8000
$ws.Cells.Item(21, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(21, 19).Value = 533
$ws.Cells.Item(21, 20).Value = 15

# Row 22
$ws.Cells.Item(22, 4).Value = 44372
$ws.Cells.Item(22, 12).Value = 'Especial'
$ws.Cells.Item(22, 13).Value = 55
$ws.Cells.Item(22, 14).Value = 15000
$ws.Cells.Item(22, 15).Value = 15000
$ws.Cells.Item(22, 16).Value = 15000
$ws.Cells.Item(22, 19).Value = 1000

# Row 23
$ws.Cells.Item(23, 4).Value = 44372
$ws.Cells.Item(23, 13).Value = 70
$ws.Cells.Item(23, 14).Value = 12000
$ws.Cells.Item(23, 15).Value = 12000
$ws.Cells.Item(23, 16).Value = 12000
$ws.Cells.Item(23, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(23, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(23, 19).Value = 800
$ws.Cells.Item(23, 20).Value = 15

# Row 24
$ws.Cells.Item(24, 4).Value = 44292
$ws.Cells.Item(24, 13).Value = 100
$ws.Cells.Item(24, 14).Value = 10500
$ws.Cells.Item(24, 15).Value = 11000
$ws.Cells.Item(24, 16).Value = 10775
$ws.Cells.Item(24, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(24, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(24, 19).Value = 599
$ws.Cells.Item(24, 20).Value = 18

# Row 25
$ws.Cells.Item(25, 4).Value = 44299
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 80

# Row 26
$ws.Cells.Item(26, 4).Value = 44316
$ws.Cells.Item(26, 13).Value = 4
$ws.Cells.Item(26, 14).Value = 230000
$ws.Cells.Item(26, 15).Value = 230000
$ws.Cells.Item(26, 16).Value = 230000
$ws.Cells.Item(26, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(26, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(26, 19).Value = 511
$ws.Cells.Item(26, 20).Value = 450

# Row 27
$ws.Cells.Item(27, 4).Value = 44313
$ws.Cells.Item(27, 12).Value = 'Primera'
$ws.Cells.Item(27, 13).Value = 6
$ws.Cells.Item(27, 14).Value = 240000
$ws.Cells.Item(27, 15).Value = 240000
$ws.Cells.Item(27, 16).Value = 240000
$ws.Cells.Item(27, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(27, 20).Value = 450

# Row 28
$ws.Cells.Item(28, 4).Value = 44405
$ws.Cells.Item(28, 12).Value = 'Primera'
$ws.Cells.Item(28, 13).Value = 50
$ws.Cells.Item(28, 14).Value = 16000
$ws.Cells.Item(28, 15).Value = 16000
$ws.Cells.Item(28, 16).Value = 16000
$ws.Cells.Item(28, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(28, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(28, 19).Value = 889
$ws.Cells.Item(28, 20).Value = 18

# Row 29
$ws.Cells.Item(29, 4).Value = 44405
$ws.Cells.Item(29, 12).Value = 'Segunda'
$ws.Cells.Item(29, 13).Value = 70
$ws.Cells.Item(29, 14).Value = 12500
$ws.Cells.Item(29, 15).Value = 12500
$ws.Cells.Item(29, 16).Value = 12500
$ws.Cells.Item(29, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(29, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(29, 19).Value = 694
$ws.Cells.Item(29, 20).Value = 18

# Row 30
$ws.Cells.Item(30, 4).Value = 44301
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 150
$ws.Cells.Item(30, 14).Value = 12000
$ws.Cells.Item(30, 15).Value = 12000
$ws.Cells.Item(30, 16).Value = 12000
$ws.Cells.Item(30, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(30, 18).Value = 'Provincia de Cachapoal'
$ws.Cells.Item(30, 19).Value = 667
$ws.Cells.Item(30, 20).Value = 18

# Row 31
$ws.Cells.Item(31, 4).Value = 44301
$ws.Cells.Item(31, 12).Value = 'Segunda'
$ws.Cells.Item(31, 14).Value = 10000
$ws.Cells.Item(31, 15).Value = 10000
$ws.Cells.Item(31, 16).Value = 10000
$ws.Cells.Item(31, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(31, 18).Value = 'Provincia de Cachapoal'
$ws.Cells.Item(31, 19).Value = 556
$ws.Cells.Item(31, 20).Value = 18

# Row 32
$ws.Cells.Item(32, 4).Value = 44327
$ws.Cells.Item(32, 12).Value = 'Primera'
$ws.Cells.Item(32, 13).Value = 4
$ws.Cells.Item(32, 14).Value = 150000
$ws.Cells.Item(32, 15).Value = 150000
$ws.Cells.Item(32, 16).Value = 150000
$ws.Cells.Item(32, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(32, 18).Value = 'Provincia de Cachapoal'
$ws.Cells.Item(32, 19).Value = 333
$ws.Cells.Item(32, 20).Value = 450

# Row 33
$ws.Cells.Item(33, 4).Value = 44384
$ws.Cells.Item(33, 12).Value = 'Especial'
$ws.Cells.Item(33, 13).Value = 70
$ws.Cells.Item(33, 14).Value = 14000
$ws.Cells.Item(33, 15).Value = 14000
$ws.Cells.Item(33, 16).Value = 14000
$ws.Cells.Item(33, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(33, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(33, 19).Value = 933
$ws.Cells.Item(33, 20).Value = 15

# Row 34
$ws.Cells.Item(34, 4).Value = 44384
$ws.Cells.Item(34, 12).Value = 'Primera'
$ws.Cells.Item(34, 13).Value = 100
$ws.Cells.Item(34, 14).Value = 12000
$ws.Cells.Item(34, 15).Value = 12000
$ws.Cells.Item(34, 16).Value = 12000
$ws.Cells.Item(34, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(34, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(34, 19).Value = 800
$ws.Cells.Item(34, 20).Value = 15

# Row 35
$ws.Cells.Item(35, 4).Value = 44384
$ws.Cells.Item(35, 12).Value = 'Segunda'
$ws.Cells.Item(35, 13).Value = 80
$ws.Cells.Item(35, 14).Value = 10000
$ws.Cells.Item(35, 15).Value = 10000
$ws.Cells.Item(35, 16).Value = 10000
$ws.Cells.Item(35, 17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(35, 19).Value = 667
$ws.Cells.Item(35, 20).Value = 15
